$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6: Week 5 summary row.
#   - G6 text: "(gained one day)" -> "(gained one day + sick for one day)"
#   - F6: new total-hours formula for the Week 5 block (B68:B78)
# ---------------------------------------------------------------------------
$ws.Range("G6").Value = "(gained one day + sick for one day)"
$ws.Range("F6").Formula = "=SUM(B68:B78)"

# ---------------------------------------------------------------------------
# Row 72: used to read
#   A72 = "Putting in game/Optimizing Audio"
#   D72 = "started at 12am"
# now becomes
#   A72 = "Putting in game/Optimizing Audio/ Fixing Neutrality"
#   B72 = 6 (hours)
#   D72 cleared (the stray "started at 12am" note is removed)
# ---------------------------------------------------------------------------
$ws.Range("A72").Value = "Putting in game/Optimizing Audio/ Fixing Neutrality"
$ws.Range("B72").Value = 6
$ws.Range("D72").ClearContents()

# ---------------------------------------------------------------------------
# New log rows 73-81: Working-on text (A), hours (B), date (C).
# Dates are written as the same serial numbers Excel already uses in this
# workbook; copying the date format from an existing dated cell (C2) first
# keeps the new cells on the workbook's existing date style instead of
# minting a duplicate number format.
# ---------------------------------------------------------------------------
$newRows = @(
  @{ Row = 73; Text = 'Fixing Neutrality/Fixing Glitches/"Touching" up game'; Hours = 2.25; Date = 41945 },
  @{ Row = 74; Text = 'Making Graphs randomly Generated'; Hours = 5; Date = 41945 },
  @{ Row = 75; Text = 'Making Graphs randomly Generated & fixing issues'; Hours = 4; Date = 41946 },
  @{ Row = 76; Text = 'fixing issues/bugs/cleaning up code'; Hours = 3; Date = 41946 },
  @{ Row = 77; Text = 'fixing issues/bugs/cleaning up code'; Hours = 5; Date = 41946 },
  @{ Row = 78; Text = 'Adding Credits + above stuff'; Hours = 3; Date = 41946 },
  @{ Row = 79; Text = 'Working on Getting in End Game stuff'; Hours = 3; Date = 41947 },
  @{ Row = 80; Text = 'Working on Getting in End Game stuff'; Hours = 9; Date = 41948 },
  @{ Row = 81; Text = 'Working on Getting in End Game stuff'; Hours = 5.5; Date = 41949 }
)

foreach ($r in $newRows) {
  $ws.Range("A" + $r.Row).Value = $r.Text
  $ws.Range("B" + $r.Row).Value = $r.Hours
  $ws.Range("C2").Copy($ws.Range("C" + $r.Row))
  $ws.Range("C" + $r.Row).Value = $r.Date
}

# ---------------------------------------------------------------------------
# Row 82: final new entry. B82 is the literal text "2+" (not a number).
# ---------------------------------------------------------------------------
$ws.Range("A82").Value = "Updating Options/Checking End Scene Works & Fixing"
$ws.Range("B82").Value = "2+"
$ws.Range("C2").Copy($ws.Range("C82"))
$ws.Range("C82").Value = 41950

# ---------------------------------------------------------------------------
# Sheet view: scroll position + active selection follow the newly added data.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 66
$ws.Range("E82").Select()
